$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$refStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = "25.778.71"
$ws.Range("E2").Value = "  -5.34%  "
$ws.Range("D3").Value = "1.812.14"
$ws.Range("E3").Value = "  -4.44%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("D4").Style = $refStyle
$ws.Range("E4").Value = "  -0.32%  "
$ws.Range("D5").Value = "'276.18"
$ws.Range("D5").Style = $refStyle
$ws.Range("E5").Value = "  -9.78%  "
$ws.Range("D6").Value = "'0.9993"
$ws.Range("D6").Style = $refStyle
$ws.Range("E6").Value = "  -0.22%  "
$ws.Range("D7").Value = "'0.5085"
$ws.Range("D7").Style = $refStyle
$ws.Range("E7").Value = "  -5.62%  "
$ws.Range("D8").Value = "'0.3520"
$ws.Range("D8").Style = $refStyle
$ws.Range("E8").Value = "  -7.20%  "
$ws.Range("B9").Value = "OKB"
$ws.Range("C9").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D9").Value = "'44.63"
$ws.Range("D9").Style = $refStyle
$ws.Range("E9").Value = "  -2.52%  "
$ws.Range("B10").Value = "Dogecoin"
$ws.Range("C10").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D10").Value = "'0.06660"
$ws.Range("D10").Style = $refStyle
$ws.Range("E10").Value = "  -8.54%  "
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D11").Value = "'20.05"
$ws.Range("D11").Style = $refStyle
$ws.Range("E11").Value = "  -8.80%  "
$ws.Range("B12").Value = "Polygon"
$ws.Range("C12").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D12").Value = "'0.8352"
$ws.Range("D12").Style = $refStyle
$ws.Range("E12").Value = "  -7.03%  "
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").Value = "'0.07814"
$ws.Range("D13").Style = $refStyle
$ws.Range("E13").Value = "  -4.57%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").Value = "1.795.37"
$ws.Range("E14").Value = "  +12.72%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").Value = "'5.072"
$ws.Range("D15").Style = $refStyle
$ws.Range("E15").Value = "  -5.12%  "
$ws.Range("B16").Value = "Litecoin"
$ws.Range("C16").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D16").Value = "'87.87"
$ws.Range("D16").Style = $refStyle
$ws.Range("E16").Value = "  -7.60%  "
$ws.Range("B17").Value = "BinanceUSD"
$ws.Range("C17").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D17").Value = "'0.9986"
$ws.Range("D17").Style = $refStyle
$ws.Range("E17").Value = "  -0.49%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").Value = "'13.88"
$ws.Range("D18").Style = $refStyle
$ws.Range("E18").Value = "  -6.50%  "
$ws.Range("D19").Value = "'0.000008021"
$ws.Range("D19").Style = $refStyle
$ws.Range("E19").Value = "  -7.12%  "
$ws.Range("B20").Value = "Dai"
$ws.Range("C20").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D20").Value = "'0.9997"
$ws.Range("D20").Style = $refStyle
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("B21").Value = "WrappedBTC"
$ws.Range("C21").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D21").Value = "25.831.59"
$ws.Range("E21").Value = "  -5.32%  "
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").Value = "'4.731"
$ws.Range("D22").Style = $refStyle
$ws.Range("E22").Value = "  -6.11%  "
$ws.Range("B23").Value = "Cosmos"
$ws.Range("C23").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D23").Value = "'9.994"
$ws.Range("D23").Style = $refStyle
$ws.Range("E23").Value = "  -7.63%  "
$ws.Range("B24").Value = "Chainlink"
$ws.Range("C24").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D24").Value = "'6.055"
$ws.Range("D24").Style = $refStyle
$ws.Range("E24").Value = "  -6.77%  "
$ws.Range("D25").Value = "'2.215"
$ws.Range("D25").Style = $refStyle
$ws.Range("E25").Value = "  -3.60%  "
$ws.Range("B26").Value = "Monero"
$ws.Range("C26").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D26").Value = "'141.81"
$ws.Range("D26").Style = $refStyle
$ws.Range("E26").Value = "  -4.57%  "
$ws.Range("B27").Value = "Toncoin"
$ws.Range("C27").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D27").Value = "'1.654"
$ws.Range("D27").Style = $refStyle
$ws.Range("E27").Value = "  -5.49%  "
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").Value = "'17.02"
$ws.Range("D28").Style = $refStyle
$ws.Range("E28").Value = "  -7.23%  "
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").Value = "'108.77"
$ws.Range("D29").Style = $refStyle
$ws.Range("E29").Value = "  -6.41%  "
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").Value = "'4.338"
$ws.Range("D30").Style = $refStyle
$ws.Range("E30").Value = "  -10.00%  "
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").Value = "'4.212"
$ws.Range("D31").Style = $refStyle
$ws.Range("E31").Value = "  -9.44%  "
$ws.Range("B32").Value = "Stellar"
$ws.Range("C32").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D32").Value = "'0.08782"
$ws.Range("D32").Style = $refStyle
$ws.Range("E32").Value = "  -4.19%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.04873"
$ws.Range("D33").Style = $refStyle
$ws.Range("E33").Value = "  -3.55%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7338"
$ws.Range("D34").Style = $refStyle
$ws.Range("E34").Value = "  -10.65%  "
$ws.Range("D35").Value = "'1.139"
$ws.Range("D35").Style = $refStyle
$ws.Range("E35").Value = "  -6.56%  "
$ws.Range("B36").Value = "HuobiToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D36").Value = "'2.895"
$ws.Range("D36").Style = $refStyle
$ws.Range("E36").Value = "  -4.19%  "
$ws.Range("B37").Value = "Frax"
$ws.Range("C37").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D37").Value = "'0.9987"
$ws.Range("D37").Style = $refStyle
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("B38").Value = "MXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D38").Value = "'3.029"
$ws.Range("D38").Style = $refStyle
$ws.Range("E38").Value = "  -8.03%  "
$ws.Range("B39").Value = "TheSandbox"
$ws.Range("C39").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D39").Value = "'0.5222"
$ws.Range("D39").Style = $refStyle
$ws.Range("E39").Value = "  -12.36%  "
$ws.Range("B40").Value = "VeChain"
$ws.Range("C40").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D40").Value = "'0.01853"
$ws.Range("D40").Style = $refStyle
$ws.Range("E40").Value = "  -6.71%  "
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").Value = "'2.298"
$ws.Range("D41").Style = $refStyle
$ws.Range("E41").Value = "  -14.01%  "
$ws.Range("B42").Value = "TrustWalletToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D42").Value = "'0.9515"
$ws.Range("D42").Style = $refStyle
$ws.Range("E42").Value = "  -11.52%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").Value = "'112.66"
$ws.Range("D43").Style = $refStyle
$ws.Range("E43").Value = "  -1.85%  "
$ws.Range("B44").Value = "FraxShare"
$ws.Range("C44").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D44").Value = "'6.189"
$ws.Range("D44").Style = $refStyle
$ws.Range("E44").Value = "  -6.79%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").Value = "'8.106"
$ws.Range("D45").Style = $refStyle
$ws.Range("E45").Value = "  -12.18%  "
$ws.Range("B46").Value = "PaxDollar"
$ws.Range("C46").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D46").Value = "'0.9989"
$ws.Range("D46").Style = $refStyle
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "'0.4561"
$ws.Range("D47").Style = $refStyle
$ws.Range("E47").Value = "  -10.41%  "
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").Value = "'0.1377"
$ws.Range("D48").Style = $refStyle
$ws.Range("E48").Value = "  -9.91%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.337"
$ws.Range("D49").Style = $refStyle
$ws.Range("E49").Value = "  -8.63%  "
$ws.Range("B50").Value = "Elrond"
$ws.Range("C50").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D50").Value = "'36.22"
$ws.Range("D50").Style = $refStyle
$ws.Range("E50").Value = "  -4.53%  "
$ws.Range("B51").Value = "NEARProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D51").Value = "'1.499"
$ws.Range("D51").Style = $refStyle
$ws.Range("E51").Value = "  -7.93%  "
